# One-click update from Van Paper 07:10 AM on 2025-11-12
#
# Inserts a new leaderboard row for customer "NICOLLET COURT RETAIL MALL"
# (Salesperson: Steiner, Owen A; Prospect code: 015; Customer Number:
# 0008368; no Last Invoice Date yet) directly above "HOLY FAMILY MARONITE
# CHURCH", pushing it and the rows below down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "HOLY FAMILY MARONITE CHURCH" currently lives in row 28 - insert a new
# blank row there (shifting it and everything below down to row 29+).
$ws.Rows("28:28").Insert()

# Match the row height used by the rest of the data rows (the new row
# otherwise inherits the sheet's default height).
$ws.Rows("28:28").RowHeight = 13.05

$ws.Cells.Item(28, 1).Value = "NICOLLET COURT RETAIL MALL"
$ws.Cells.Item(28, 2).Value = "Steiner, Owen A"
$ws.Cells.Item(28, 3).Value = "015"
$ws.Cells.Item(28, 5).Value = "0008368"
